$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.021.66'
$ws.Cells.Item(2, 5).Value = '  -0.37%  '
$ws.Cells.Item(3, 4).Value = '2.551.57'
$ws.Cells.Item(3, 5).Value = '  +0.23%  '
$ws.Cells.Item(4, 5).Value = '  -0.12%  '
$ws.Cells.Item(5, 4).Value = '583.04'
$ws.Cells.Item(5, 5).Value = '  +2.33%  '
$ws.Cells.Item(6, 4).Value = '147.51'
$ws.Cells.Item(6, 5).Value = '  -1.88%  '
$ws.Cells.Item(7, 5).Value = '  -0.07%  '
$ws.Cells.Item(8, 5).Value = '  -0.91%  '
$ws.Cells.Item(9, 5).Value = '  -0.25%  '
$ws.Cells.Item(10, 5).Value = '  -3.43%  '
$ws.Cells.Item(11, 5).Value = '  -0.02%  '
$ws.Cells.Item(12, 5).Value = '  -0.84%  '
$ws.Cells.Item(13, 4).Value = '27.26'
$ws.Cells.Item(13, 5).Value = '  -3.53%  '
$ws.Cells.Item(14, 4).Value = '3.007.46'
$ws.Cells.Item(14, 5).Value = '  +0.16%  '
$ws.Cells.Item(15, 4).Value = '62.935.37'
$ws.Cells.Item(15, 5).Value = '  -0.48%  '
$ws.Cells.Item(16, 5).Value = '  -0.13%  '
$ws.Cells.Item(17, 4).Value = '2.552.89'
$ws.Cells.Item(17, 5).Value = '  +0.30%  '
$ws.Cells.Item(18, 4).Value = '11.35'
$ws.Cells.Item(18, 5).Value = '  -2.20%  '
$ws.Cells.Item(19, 4).Value = '336.21'
$ws.Cells.Item(19, 5).Value = '  -1.39%  '
$ws.Cells.Item(20, 4).Value = '4.34'
$ws.Cells.Item(20, 5).Value = '  -0.55%  '
$ws.Cells.Item(21, 5).Value = '  -1.41%  '
$ws.Cells.Item(22, 5).Value = '  +0.09%  '
$ws.Cells.Item(23, 5).Value = '  -0.22%  '
$ws.Cells.Item(24, 2).Value = 'Kaspa'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(24, 4).Value = '0.170'
$ws.Cells.Item(24, 5).Value = '  -0.32%  '
$ws.Cells.Item(25, 2).Value = 'Fetch.AI'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(25, 4).Value = '1.63'
$ws.Cells.Item(25, 5).Value = '  +1.75%  '
$ws.Cells.Item(26, 2).Value = 'SuiNetwork'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(26, 4).Value = '1.49'
$ws.Cells.Item(26, 5).Value = '  +1.38%  '
$ws.Cells.Item(27, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(27, 4).Value = '1.00'
$ws.Cells.Item(27, 5).Value = '  -0.08%  '
$ws.Cells.Item(28, 4).Value = '8.38'
$ws.Cells.Item(28, 5).Value = '  -0.68%  '
$ws.Cells.Item(29, 4).Value = '7.42'
$ws.Cells.Item(29, 5).Value = '  +4.45%  '
$ws.Cells.Item(30, 5).Value = '  +2.81%  '
$ws.Cells.Item(31, 4).Value = '0.0₃0815'
$ws.Cells.Item(31, 5).Value = '  -1.75%  '
$ws.Cells.Item(32, 4).Value = '177.90'
$ws.Cells.Item(32, 5).Value = '  +1.08%  '
$ws.Cells.Item(33, 5).Value = '  -0.78%  '
$ws.Cells.Item(34, 4).Value = '409.48'
$ws.Cells.Item(34, 5).Value = '  -2.20%  '
$ws.Cells.Item(35, 2).Value = 'EthereumClassic'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(35, 4).Value = '19.16'
$ws.Cells.Item(35, 5).Value = '  +0.32%  '
$ws.Cells.Item(36, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(36, 4).Value = '0.400'
$ws.Cells.Item(36, 5).Value = '  -1.33%  '
$ws.Cells.Item(37, 5).Value = '  +0.02%  '
$ws.Cells.Item(38, 5).Value = '  -1.55%  '
$ws.Cells.Item(39, 5).Value = '  -0.19%  '
$ws.Cells.Item(40, 5).Value = '  -0.11%  '
$ws.Cells.Item(41, 4).Value = '39.60'
$ws.Cells.Item(41, 5).Value = '  -1.12%  '
$ws.Cells.Item(42, 4).Value = '150.94'
$ws.Cells.Item(42, 5).Value = '  -2.66%  '
$ws.Cells.Item(43, 5).Value = '  -0.92%  '
$ws.Cells.Item(44, 4).Value = '20.90'
$ws.Cells.Item(44, 5).Value = '  -0.74%  '
$ws.Cells.Item(45, 4).Value = '0.0544'
$ws.Cells.Item(45, 5).Value = '  +2.52%  '
$ws.Cells.Item(46, 4).Value = '0.602'
$ws.Cells.Item(46, 5).Value = '  -0.96%  '
$ws.Cells.Item(47, 4).Value = '0.0970'
$ws.Cells.Item(47, 5).Value = '  +0.51%  '
$ws.Cells.Item(48, 4).Value = '0.0240'
$ws.Cells.Item(48, 5).Value = '  +1.55%  '
$ws.Cells.Item(49, 4).Value = '18.32'
$ws.Cells.Item(49, 5).Value = '  -1.88%  '
$ws.Cells.Item(50, 5).Value = '  -4.95%  '
$ws.Cells.Item(51, 5).Value = '  -0.06%  '
